$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Year" column (D) is stored as text in the original workbook. Writing a
# plain numeric-looking string via .Value lets Excel auto-convert it to a
# number, so we force a text number format, set the value, then restore the
# cell's original style to avoid leaving a stray "Text" style behind.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("B2").Value = "Accountability and cyber conflict: examining institutional constraints on the use of cyber proxies"
$ws.Range("C2").Value = "William Akoto"
Set-TextValue $ws.Range("D2") "2022"
$ws.Range("E2").Value = "10.1177/07388942211051264"
$ws.Range("F2").Value = "Restricted"

$ws.Range("B3").Value = "Mapping Global Cyberterror Networks: An Empirical Study of Al-Qaeda and ISIS Cyberterrorism Events"
$ws.Range("C3").Value = "Claire Seungeun Lee, Kyung-Shick Choi, Ryan Shandler, Chris Kayser"
Set-TextValue $ws.Range("D3") "2021"
$ws.Range("E3").Value = "10.1177/10439862211001606"

$ws.Range("B4").Value = "Cyber and contentious politics: Evidence from the US radical environmental movement"
$ws.Range("C4").Value = "Thomas Zeitzoff, Grace Gold"
Set-TextValue $ws.Range("D4") "2024"
$ws.Range("E4").Value = "10.1177/00223433231221426"

$ws.Range("B5").Value = "Cyberattacks and public opinion – The effect of uncertainty in guiding preferences"
$ws.Range("C5").Value = "Eric Jardine, Nathaniel Porter, Ryan Shandler"
Set-TextValue $ws.Range("D5") "2024"
$ws.Range("E5").Value = "10.1177/00223433231218178"

$ws.Range("B6").Value = "A virtual necessity: Some modest steps toward greater cybersecurity"
$ws.Range("C6").Value = "Herbert Lin"
Set-TextValue $ws.Range("D6") "2012"
$ws.Range("E6").Value = "10.1177/0096340212459039"

$ws.Range("B7").Value = "Global versus Local Optimization in Redundancy Resolution of Robotic Manipulators"
$ws.Range("C7").Value = "Kazem Kazerounian, Zhaoyu Wang"
Set-TextValue $ws.Range("D7") "1988"
$ws.Range("E7").Value = "10.1177/027836498800700501"

$ws.Range("B8").Value = "From the Ontology of Video Games to the Epistemology of Digital Movements. Towards a Semiotics of Virtual Practices"
$ws.Range("C8").Value = "Enzo D’Armenio"
Set-TextValue $ws.Range("D8") "2024"
$ws.Range("E8").Value = "10.1177/15554120241263630"
$ws.Range("F8").Value = "Restricted"

$ws.Range("B9").Value = "Tech titans, cyber commons and the war in Ukraine: An incipient shift in international relations"
$ws.Range("C9").Value = "Eviatar Matania, Udi Sommer"
Set-TextValue $ws.Range("D9") "2023"
$ws.Range("E9").Value = "10.1177/00471178231211500"
$ws.Range("F9").Value = "Open Access"

$ws.Range("B10").Value = "Warring from the virtual to the real: Assessing the public’s threshold for war over cyber security"
$ws.Range("C10").Value = "Sarah Kreps, Debak Das"
Set-TextValue $ws.Range("D10") "2017"
$ws.Range("E10").Value = "10.1177/2053168017715930"
$ws.Range("F10").Value = "Open Access"

$ws.Range("B11").Value = "Cyber-Flirting: Playing at Love on the Internet"
$ws.Range("C11").Value = "Monica Therese Whitty"
Set-TextValue $ws.Range("D11") "2003"
$ws.Range("E11").Value = "10.1177/0959354303013003003"
